# "Order table implement and methods added"
# Builds out the Orders table (rows 31-43) with its columns, API entry
# points/methods, and hyperlinks, mirroring the structure already used by
# the User/Cart/Product tables above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 31: orderId -------------------------------------------------
$ws.Range("B31").Value = "orderId"

# E31 did not exist before; copy the plain "SITE column" formatting (style
# used throughout column E, e.g. E34) before writing into it.
$ws.Range("E34").Copy()
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("E31").Value = "Cart Table"

$ws.Range("F31").Value = "Submit a order"
$ws.Hyperlinks.Add($ws.Range("G31"), "http://localhost:3005/orders/submitNew")

# --- Row 32: Product Name / SITE / User Table / Get all the orders ---
$ws.Range("B32").Value = "Product Name"
$ws.Range("D32").Value = "SITE"

# E32 keeps its existing text-only look (no border/fill) - apply the same
# number format / font as the rest of the table but strip border & fill.
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Font.Size = 12
$ws.Range("E32").Font.Bold = $false
$ws.Range("E32").Borders.LineStyle = -4142
$ws.Range("E32").Interior.Pattern = -4142
$ws.Range("E32").Value = "User Table"

$ws.Range("F32").Value = "Get all the orders"
$ws.Hyperlinks.Add($ws.Range("G32"), "http://localhost:3005/orders")

# H32 did not exist before; give it the same look as H31/H33's column
# (method cells use style from H31).
$ws.Range("H31").Copy()
$ws.Range("H32").PasteSpecial(-4122)
$ws.Range("H32").Value = "GET"

# --- Row 33: Product Image --------------------------------------------
$ws.Range("B33").Value = "Product Image"

# F33:H33 did not exist; pick up the blank "SITE row" formatting used by
# the rest of the table (e.g. D34/E34 style == s4).
$ws.Range("D34").Copy()
$ws.Range("F33:H33").PasteSpecial(-4122)

# --- Row 35: UserID ------------------------------------------------------
$ws.Range("B35").Value = "UserID"

# --- New rows 36-43: EmailID, User Name, Phone, Addressline, Land Mark,
#     Zip, City, State (mirrors the User table's sub-fields) -------------
$ws.Range("A34:H34").Copy()
$ws.Range("A36:H43").PasteSpecial(-4122)

$ws.Range("B36").Value = "EmailID"
$ws.Range("D36").Value = "SITE"

$ws.Range("B37").Value = "User Name"
$ws.Range("D37").Value = "SITE"

$ws.Range("B38").Value = "Phone"
$ws.Range("D38").Value = "SITE"

$ws.Range("B39").Value = "Addressline"
$ws.Range("D39").Value = "SITE"

$ws.Range("B40").Value = "Land Mark"
$ws.Range("D40").Value = "SITE"

$ws.Range("B41").Value = "Zip"
$ws.Range("D41").Value = "SITE"

$ws.Range("B42").Value = "City"
$ws.Range("D42").Value = "SITE"

$ws.Range("B43").Value = "State"
$ws.Range("D43").Value = "SITE"

# --- View state: leave the cursor parked near the bottom of the new table
$ws.Range("H38").Select()
